# Auto-generated Excel COM-interop script to refresh Market Board price columns
# (currentAveragePrice / currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 795.5789
$ws.Range("I28").Value = 707.25
$ws.Range("J28").Value = 1266.6666
$ws.Range("K28").Value = 707.25
$ws.Range("L28").Value = 1266.6666
$ws.Range("M28").Value = -222.25
$ws.Range("N28").Value = -2236.6666

$ws.Range("H40").Value = 1161.5385
$ws.Range("J40").Value = 1161.5385
$ws.Range("L40").Value = 1161.5385
$ws.Range("N40").Value = -1511.5385

$ws.Range("H64").Value = 3609.4285
$ws.Range("I64").Value = 3504.875
$ws.Range("J64").Value = 3697.4736
$ws.Range("K64").Value = 3504.875
$ws.Range("L64").Value = 3697.4736
$ws.Range("M64").Value = -3256.875
$ws.Range("N64").Value = -4193.473599999999

$ws.Range("H67").Value = 3609.4285
$ws.Range("I67").Value = 3504.875
$ws.Range("J67").Value = 3697.4736
$ws.Range("K67").Value = 3504.875
$ws.Range("L67").Value = 3697.4736
$ws.Range("M67").Value = -2646.875
$ws.Range("N67").Value = -5413.473599999999

$ws.Range("H76").Value = 7582.143
$ws.Range("I76").Value = 8144.7915
$ws.Range("K76").Value = 8144.7915
$ws.Range("M76").Value = -7829.7915

$ws.Range("H79").Value = 7582.143
$ws.Range("I79").Value = 8144.7915
$ws.Range("K79").Value = 8144.7915
$ws.Range("M79").Value = -7052.7915

$ws.Range("H112").Value = 1894.3103
$ws.Range("I112").Value = 1358.75
$ws.Range("J112").Value = 1980
$ws.Range("K112").Value = 4076.25
$ws.Range("L112").Value = 5940
$ws.Range("M112").Value = -2968.25
$ws.Range("N112").Value = -8156

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 39936
$ws.Range("I9").Value = 38904
$ws.Range("K9").Value = 38904
$ws.Range("M9").Value = -38734

$ws.Range("H20").Value = 39936
$ws.Range("I20").Value = 38904
$ws.Range("K20").Value = 38904
$ws.Range("M20").Value = -38634

$ws.Range("H32").Value = 3122.06
$ws.Range("I32").Value = 2682.3447
$ws.Range("K32").Value = 2682.3447
$ws.Range("M32").Value = -2395.3447

$ws.Range("H74").Value = 1812.6852
$ws.Range("I74").Value = 1087.925
$ws.Range("J74").Value = 3883.4285
$ws.Range("K74").Value = 1087.925
$ws.Range("L74").Value = 3883.4285
$ws.Range("M74").Value = -213.925
$ws.Range("N74").Value = -5631.4285

$ws.Range("H77").Value = 1812.6852
$ws.Range("I77").Value = 1087.925
$ws.Range("J77").Value = 3883.4285
$ws.Range("K77").Value = 5439.625
$ws.Range("L77").Value = 19417.1425
$ws.Range("M77").Value = -1071.625
$ws.Range("N77").Value = -28153.1425

$ws.Range("H110").Value = 1562.1852
$ws.Range("I110").Value = 1642.9412
$ws.Range("J110").Value = 1424.9
$ws.Range("K110").Value = 1642.9412
$ws.Range("L110").Value = 1424.9
$ws.Range("M110").Value = 402.0588
$ws.Range("N110").Value = -5514.9

$ws.Range("H132").Value = 3870.077
$ws.Range("I132").Value = 3947
$ws.Range("J132").Value = 3759.5
$ws.Range("K132").Value = 11841
$ws.Range("L132").Value = 11278.5
$ws.Range("M132").Value = -9311
$ws.Range("N132").Value = -16338.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 365.08694
$ws.Range("I80").Value = 253.125
$ws.Range("J80").Value = 424.8
$ws.Range("K80").Value = 253.125
$ws.Range("L80").Value = 424.8
$ws.Range("M80").Value = 744.875
$ws.Range("N80").Value = -2420.8

$ws.Range("H83").Value = 365.08694
$ws.Range("I83").Value = 253.125
$ws.Range("J83").Value = 424.8
$ws.Range("K83").Value = 1265.625
$ws.Range("L83").Value = 2124
$ws.Range("M83").Value = 3726.375
$ws.Range("N83").Value = -12108

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3409.5925
$ws.Range("I31").Value = 2487.9285
$ws.Range("J31").Value = 4402.154
$ws.Range("K31").Value = 2487.9285
$ws.Range("L31").Value = 4402.154
$ws.Range("M31").Value = -2192.9285
$ws.Range("N31").Value = -4992.154

$ws.Range("H34").Value = 3409.5925
$ws.Range("I34").Value = 2487.9285
$ws.Range("J34").Value = 4402.154
$ws.Range("K34").Value = 2487.9285
$ws.Range("L34").Value = 4402.154
$ws.Range("M34").Value = -2285.9285
$ws.Range("N34").Value = -4806.154

$ws.Range("H39").Value = 7500
$ws.Range("I39").Value = 7500
$ws.Range("K39").Value = 7500
$ws.Range("M39").Value = -7109

$ws.Range("H49").Value = 7500
$ws.Range("I49").Value = 7500
$ws.Range("K49").Value = 7500
$ws.Range("M49").Value = -7318

$ws.Range("H62").Value = 3199.8223
$ws.Range("I62").Value = 2990.5908
$ws.Range("J62").Value = 3399.9565
$ws.Range("K62").Value = 2990.5908
$ws.Range("L62").Value = 3399.9565
$ws.Range("M62").Value = -2366.5908
$ws.Range("N62").Value = -4647.9565

$ws.Range("H65").Value = 3199.8223
$ws.Range("I65").Value = 2990.5908
$ws.Range("J65").Value = 3399.9565
$ws.Range("K65").Value = 14952.954
$ws.Range("L65").Value = 16999.7825
$ws.Range("M65").Value = -11832.954
$ws.Range("N65").Value = -23239.7825

$ws.Range("H105").Value = 619
$ws.Range("I105").Value = 605.9091
$ws.Range("J105").Value = 715
$ws.Range("K105").Value = 605.9091
$ws.Range("L105").Value = 715
$ws.Range("M105").Value = 1141.0909
$ws.Range("N105").Value = -4209

$ws.Range("H132").Value = 2030.7778
$ws.Range("I132").Value = 937.8214
$ws.Range("J132").Value = 3830.9412
$ws.Range("K132").Value = 2813.4642
$ws.Range("L132").Value = 11492.8236
$ws.Range("M132").Value = -283.4642000000003
$ws.Range("N132").Value = -16552.8236

$ws.Range("H134").Value = 1654.7291
$ws.Range("I134").Value = 980.069
$ws.Range("J134").Value = 2684.4736
$ws.Range("K134").Value = 2940.207
$ws.Range("L134").Value = 8053.4208
$ws.Range("M134").Value = -405.2069999999999
$ws.Range("N134").Value = -13123.4208

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 970.875
$ws.Range("J34").Value = 1027.3182
$ws.Range("L34").Value = 3081.9546
$ws.Range("N34").Value = -3249.9546

$ws.Range("H56").Value = 6200
$ws.Range("I56").Value = 6200
$ws.Range("K56").Value = 6200
$ws.Range("M56").Value = -5670

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 533852
$ws.Range("I14").Value = 533852
$ws.Range("K14").Value = 533852
$ws.Range("M14").Value = -533684

$ws.Range("H70").Value = 7433.394
$ws.Range("I70").Value = 3824.12
$ws.Range("J70").Value = 18712.375
$ws.Range("K70").Value = 3824.12
$ws.Range("L70").Value = 18712.375
$ws.Range("M70").Value = -3554.12
$ws.Range("N70").Value = -19252.375

$ws.Range("H73").Value = 7433.394
$ws.Range("I73").Value = 3824.12
$ws.Range("J73").Value = 18712.375
$ws.Range("K73").Value = 3824.12
$ws.Range("L73").Value = 18712.375
$ws.Range("M73").Value = -2888.12
$ws.Range("N73").Value = -20584.375

$ws.Range("H80").Value = 4281.116
$ws.Range("I80").Value = 4529.4443
$ws.Range("J80").Value = 3004
$ws.Range("K80").Value = 4529.4443
$ws.Range("L80").Value = 3004
$ws.Range("M80").Value = -3531.4443
$ws.Range("N80").Value = -5000

$ws.Range("H83").Value = 4281.116
$ws.Range("I83").Value = 4529.4443
$ws.Range("J83").Value = 3004
$ws.Range("K83").Value = 22647.2215
$ws.Range("L83").Value = 15020
$ws.Range("M83").Value = -17655.2215
$ws.Range("N83").Value = -25004

$ws.Range("H126").Value = 2885.1333
$ws.Range("I126").Value = 2762.6365
$ws.Range("J126").Value = 3222
$ws.Range("K126").Value = 8287.9095
$ws.Range("L126").Value = 9666
$ws.Range("M126").Value = -5817.9095
$ws.Range("N126").Value = -14606

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 2615
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws.Range("H28").Value = 2615
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()

$ws.Range("H30").Value = 36
$ws.Range("I30").Value = 36
$ws.Range("K30").Value = 36
$ws.Range("M30").Value = 72

$ws.Range("H37").Value = 2615
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()

$ws.Range("H43").Value = 9300
$ws.Range("I43").Value = 8600
$ws.Range("J43").Value = 10000
$ws.Range("K43").Value = 8600
$ws.Range("L43").Value = 10000
$ws.Range("M43").Value = -8407
$ws.Range("N43").Value = -10386

$ws.Range("H46").Value = 769.62964
$ws.Range("I46").Value = 629.2308
$ws.Range("J46").Value = 900
$ws.Range("K46").Value = 629.2308
$ws.Range("L46").Value = 900
$ws.Range("M46").Value = -441.2308
$ws.Range("N46").Value = -1276

$ws.Range("H110").Value = 12455.333
$ws.Range("J110").Value = 12455.333
$ws.Range("L110").Value = 12455.333
$ws.Range("N110").Value = -20635.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 696.6667
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 545
$ws.Range("K107").Value = 3000
$ws.Range("L107").Value = 1635
$ws.Range("M107").Value = -1080
$ws.Range("N107").Value = -5475
